$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(-40.18527779212365, -12.29963522690993, -2.489131450653076),
    @(-5.963009345225885, -38.13061856000844, 21.80784027393036),
    @(6.794984866411211, -17.97448307428637, 8.985512146583247),
    @(8.772580513587364, -2.704863694997934, -3.532207745772141),
    @(15.409985371125, -21.02941689124464, -3.300734495505285),
    @(29.99191743899604, -16.15106448149102, -5.187779292082201),
    @(43.31018800001878, -13.67232219989484, -2.019487449756035),
    @(47.68574269612633, -58.8819274902343, 7.092501203219018),
    @(27.69281722337806, -49.70874527173221, 14.84765016115646),
    @(-3.377501120934119, 4.849835469172549, 0.1248559951782222)
)

$startRow = 22
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
